# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner price/profit updates to the Asura_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 795.1
$ws.Range("I6").Value = 562.75
$ws.Range("K6").Value = 1688.25
$ws.Range("M6").Value = -1576.25

$ws.Range("H11").Value = 49.375
$ws.Range("I11").Value = 49.375
$ws.Range("K11").Value = 49.375
$ws.Range("M11").Value = 90.625

$ws.Range("H31").Value = 1445763.4
$ws.Range("I31").Value = 2020268.8
$ws.Range("J31").Value = 9500
$ws.Range("K31").Value = 6060806.4
$ws.Range("L31").Value = 28500
$ws.Range("M31").Value = -6060576.4
$ws.Range("N31").Value = -28960

$ws.Range("H38").Value = 995
$ws.Range("I38").Value = 75
$ws.Range("J38").Value = 3525
$ws.Range("K38").Value = 225
$ws.Range("L38").Value = 10575
$ws.Range("M38").Value = 147
$ws.Range("N38").Value = -11319

$ws.Range("H39").Value = 213.6
$ws.Range("I39").Value = 217
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 651
$ws.Range("L39").Value = 600
$ws.Range("M39").Value = -355
$ws.Range("N39").Value = -1192

$ws.Range("H70").Value = 67973.664
$ws.Range("I70").Value = 251075
$ws.Range("J70").Value = 1391.3636
$ws.Range("K70").Value = 753225
$ws.Range("L70").Value = 4174.0908
$ws.Range("M70").Value = -752955
$ws.Range("N70").Value = -4714.0908

$ws.Range("H73").Value = 67973.664
$ws.Range("I73").Value = 251075
$ws.Range("J73").Value = 1391.3636
$ws.Range("K73").Value = 753225
$ws.Range("L73").Value = 4174.0908
$ws.Range("M73").Value = -752289
$ws.Range("N73").Value = -6046.0908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1216.6666
$ws.Range("I2").Value = 1050
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 1050
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = -937
$ws.Range("N2").Value = -1526

$ws.Range("H32").Value = 17083.691
$ws.Range("I32").Value = 15907.63
$ws.Range("J32").Value = 22857.092
$ws.Range("K32").Value = 15907.63
$ws.Range("L32").Value = 22857.092
$ws.Range("M32").Value = -15620.63
$ws.Range("N32").Value = -23431.092

$ws.Range("H74").Value = 1515.9062
$ws.Range("I74").Value = 1306.3334
$ws.Range("J74").Value = 1785.3572
$ws.Range("K74").Value = 1306.3334
$ws.Range("L74").Value = 1785.3572
$ws.Range("M74").Value = -432.3334
$ws.Range("N74").Value = -3533.3572

$ws.Range("H77").Value = 1515.9062
$ws.Range("I77").Value = 1306.3334
$ws.Range("J77").Value = 1785.3572
$ws.Range("K77").Value = 6531.666999999999
$ws.Range("L77").Value = 8926.786
$ws.Range("M77").Value = -2163.666999999999
$ws.Range("N77").Value = -17662.786

$ws.Range("H116").Value = 1216.6666
$ws.Range("I116").Value = 1050
$ws.Range("J116").Value = 1300
$ws.Range("K116").Value = 1050
$ws.Range("L116").Value = 1300
$ws.Range("M116").Value = 1244
$ws.Range("N116").Value = -5888

$ws.Range("H132").Value = 1820185.9
$ws.Range("I132").Value = 1820185.9
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5460557.699999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5458027.699999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1216.6666
$ws.Range("I3").Value = 1050
$ws.Range("J3").Value = 1300
$ws.Range("K3").Value = 1050
$ws.Range("L3").Value = 1300
$ws.Range("M3").Value = -936
$ws.Range("N3").Value = -1528

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2736.7058
$ws.Range("I31").Value = 1903.8918
$ws.Range("K31").Value = 1903.8918
$ws.Range("M31").Value = -1608.8918

$ws.Range("H34").Value = 2736.7058
$ws.Range("I34").Value = 1903.8918
$ws.Range("K34").Value = 1903.8918
$ws.Range("M34").Value = -1701.8918

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 408.8
$ws.Range("I6").Value = 184.14285
$ws.Range("J6").Value = 933
$ws.Range("K6").Value = 552.4285500000001
$ws.Range("L6").Value = 2799
$ws.Range("M6").Value = -439.4285500000001
$ws.Range("N6").Value = -3025

$ws.Range("H131").Value = 12661022
$ws.Range("I131").Value = 9496.583000000001
$ws.Range("J131").Value = 14926967
$ws.Range("K131").Value = 28489.749
$ws.Range("L131").Value = 44780901
$ws.Range("M131").Value = -23449.749
$ws.Range("N131").Value = -44790981

$ws.Range("H134").Value = 3653.2903
$ws.Range("J134").Value = 7274.778
$ws.Range("L134").Value = 21824.334
$ws.Range("N134").Value = -31964.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 8002733
$ws.Range("J3").Value = 36673664
$ws.Range("L3").Value = 36673664
$ws.Range("N3").Value = -36673896

$ws.Range("H10").Value = 10150
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 10150
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 10150
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -10488

$ws.Range("H101").Value = 54279.5
$ws.Range("J101").Value = 54279.5
$ws.Range("L101").Value = 54279.5
$ws.Range("N101").Value = -60769.5

$ws.Range("H109").Value = 9280
$ws.Range("J109").Value = 9280
$ws.Range("L109").Value = 9280
$ws.Range("N109").Value = -11360

$ws.Range("H113").Value = 2352.625
$ws.Range("J113").Value = 2521.6
$ws.Range("L113").Value = 2521.6
$ws.Range("N113").Value = -6861.6

$ws.Range("H132").Value = 1851.5555
$ws.Range("I132").Value = 1582.125
$ws.Range("J132").Value = 4007
$ws.Range("K132").Value = 4746.375
$ws.Range("L132").Value = 12021
$ws.Range("M132").Value = -2216.375
$ws.Range("N132").Value = -17081

$ws.Range("H136").Value = 25498.357
$ws.Range("J136").Value = 25498.357
$ws.Range("L136").Value = 76495.071
$ws.Range("N136").Value = -81595.071

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 320.66666
$ws.Range("I22").Value = 360.8889
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 360.8889
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -65.88889999999998
$ws.Range("N22").Value = -790

$ws.Range("H27").Value = 320.66666
$ws.Range("I27").Value = 360.8889
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 360.8889
$ws.Range("L27").Value = 200
$ws.Range("M27").Value = -253.8889
$ws.Range("N27").Value = -414

$ws.Range("H46").Value = 1727.1818
$ws.Range("I46").Value = 1499.8334
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1499.8334
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1311.8334
$ws.Range("N46").Value = -2376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 20000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -26240

$ws.Range("H132").Value = 1049.4894
$ws.Range("I132").Value = 1011.5952
$ws.Range("J132").Value = 1367.8
$ws.Range("K132").Value = 3034.7856
$ws.Range("L132").Value = 4103.4
$ws.Range("M132").Value = -504.7856000000002
$ws.Range("N132").Value = -9163.4
